$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the data contained in row 2 and row 3 (columns A:L).
# Using Range.Value directly would cause Excel to "smart convert" numeric-
# looking / date-looking text (e.g. "7836", "1.1", "06/09/2025") into real
# numbers / dates, which would not match the original text cells. Copy +
# PasteSpecial (values only) preserves the original text representation.

$xlPasteValues = -4163

# Stage row 2 and row 3 contents in a scratch area (columns N:Y) so that
# the swap does not overwrite source data before it has been read.
$ws.Range("A2:L2").Copy()
$ws.Range("N2").PasteSpecial($xlPasteValues)

$ws.Range("A3:L3").Copy()
$ws.Range("N3").PasteSpecial($xlPasteValues)

# Clear the original rows completely (including cells that are blank, e.g.
# C3/H2/H3/L2/L3) before pasting the swapped values back in.
$ws.Range("A2:L3").ClearContents()

# Row 2 gets what used to be in row 3, and vice versa.
$ws.Range("N3:Y3").Copy()
$ws.Range("A2").PasteSpecial($xlPasteValues)

$ws.Range("N2:Y2").Copy()
$ws.Range("A3").PasteSpecial($xlPasteValues)

# Remove the scratch/staging data.
$ws.Range("N2:Y3").ClearContents()

$excel.CutCopyMode = 0
